# Update the "Comp controls" sheet so that each fluorochrome/detector
# reference picks up the "-A" (area) channel suffix, e.g. "B515" -> "B515-A",
# "G780" -> "G780-A", etc. This affects both the FCS file name column (A)
# and the Marker:Detector column (B) for rows 4-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

$values = @(
    @(4,  "Compensation Controls_B515-A Stained Control.fcs",             "LIVE GREEN:B515-A"),
    @(5,  "Compensation Controls_G560-A Stained Control.fcs",              "CD197:G560-A"),
    @(6,  "Compensation Controls_B710-A Stained Control.fcs",              "CD4:B710-A"),
    @(7,  "Compensation Controls_CD45RA PE-Cy7 G780-A Stained Control.fcs","CD45RA:G780-A"),
    @(8,  "Compensation Controls_CCR4 PE-Cy7 G780-A Stained Control.fcs",  "CD194:G780-A"),
    @(9,  "Compensation Controls_CD27 PE-Cy7 G780-A Stained Control.fcs",  "CD27:G780-A"),
    @(10, "Compensation Controls_CD11c PE-Cy7 G780-A Stained Control.fcs", "CD11c:G780-A"),
    @(11, "Compensation Controls_CCR6 PE-Cy7 G780-A Stained Control.fcs",  "CD196:G780-A"),
    @(12, "Compensation Controls_CD38 APC R660-A Stained Control.fcs",     "CD38:R660-A"),
    @(13, "Compensation Controls_CD127 AF647 R660-A Stained Control.fcs",  "CD127:R660-A"),
    @(14, "Compensation Controls_CD8 APC-H7 R780-A Stained Control.fcs",   "CD8:R780-A"),
    @(15, "Compensation Controls_CD45RO APC-H7 R780-A Stained Control.fcs","CD45RO:R780-A"),
    @(16, "Compensation Controls_CD20 APC-H7 R780-A Stained Control.fcs",  "CD20:R780-A"),
    @(17, "Compensation Controls_Lineage APC-H7 R780-A Stained Control.fcs","CD3+19+20:R780-A"),
    @(18, "Compensation Controls_V450-A Stained Control.fcs",              "CD3:V450-A"),
    @(19, "Compensation Controls_V545-A Stained Control.fcs",              "HLA-DR:V545-A")
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[2]
}
foreach ($row in $values) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
}

$ws.Range("A4").Select()
